# Finalized definition of virtual environment
# Updates the calibration results on the AR, SETAR, GARCH, TARCH and
# AR_TARCH sheets: both the numeric parameter values and the
# "abs_epsi_autocorr" autocorrelation-array text for each model.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.04562550251115292
$ws.Range("B3").Value = 0.8120396882278178
$ws.Range("B4").Value = 96.43207676280547
$ws.Range("B5").Value = '[1.0, 0.19201877694542827, 0.1983576640446433, 0.15699605721127627, 0.11437189584417722, 0.35345911830904025, 0.14534986961001894, 0.15549198717630752, 0.12305655074868245, 0.10258067720422014, 0.07826776012757741, 0.09524037177642192, 0.09740928089664322, 0.10280408920347854, 0.1495242878911461, 0.12322569587886337, 0.09317324009323355, 0.07167131479090139, 0.08380183271437734, 0.11472346774879641]'

# ---------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 0.9031296572280179
$ws.Range("B4").Value = -3.166948928744786
$ws.Range("B5").Value = 0.7502244064012349
$ws.Range("B6").Value = 69.8187561924296
$ws.Range("B7").Value = 3.836364800123929
$ws.Range("B8").Value = 0.6781113055542988
$ws.Range("B9").Value = 58.0956279966797
$ws.Range("B10").Value = '[1.0, 0.22633845846720407, 0.20083941064244692, 0.15013174284200512, 0.1680052779909548, 0.2377359815138933, 0.12632201858329176, 0.10983530541288955, 0.1063880697760478, 0.08404494016379159, 0.06842959644765147, 0.07911280419673475, 0.06454201349199791, 0.08271630948648453, 0.11885947917866424, 0.08122060444046485, 0.08135675221228136, 0.05296040179494797, 0.09230361516332905, 0.08954449453675926]'

# ---------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.02473144370188488
$ws.Range("B3").Value = 1.805247805946451
$ws.Range("B4").Value = 0.08399015283849863
$ws.Range("B5").Value = 0.9015039943398633
$ws.Range("B6").Value = '[1.0, -0.012570249480646832, -0.01851845280122588, -0.05910702294046567, -0.08551980785327967, 0.282896005165324, 0.011070542747553013, 0.03119874116986023, -0.024121303626147016, -0.017087644801969323, -0.06098195323044896, -0.051876205104514836, -0.018426770103801258, -0.00585158267805371, 0.060877216591329206, 0.012052561264814535, -0.029108120343660313, -0.05409705432783527, -0.04456293021433298, -0.0063834955291892205]'

# ---------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.004109738811114168
$ws.Range("B3").Value = 1.801905728858356
$ws.Range("B4").Value = 0.07791483242018225
$ws.Range("B5").Value = 0.9019065457830313
$ws.Range("B6").Value = '[0.9999999999999999, -0.012216490567329416, -0.018799691257434185, -0.059001843732208244, -0.08519837350833613, 0.28257030624722046, 0.011284636334402198, 0.031077568909843284, -0.02390257209279314, -0.01708269363288374, -0.06068968557000507, -0.05177131967106617, -0.018547384796799765, -0.006111562350623109, 0.061014234442891564, 0.012183700908471285, -0.02909016560562657, -0.05374892125045072, -0.04447429651496943, -0.006190288521372873]'
$ws.Range("B7").Value = 0.01124457692643716

# ---------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.1935749980606548
$ws.Range("B3").Value = 1.559249988892618
$ws.Range("B4").Value = 0.1050183655641683
$ws.Range("B5").Value = 0.9012108961626137
$ws.Range("B6").Value = '[0.9999999999999998, 0.010701332566368441, 0.025242362120509523, -0.018294170579458494, -0.0663176070816311, 0.19670832254093507, 0.0010840124992194248, 0.002668947228981423, -0.017892847807274202, -0.02644926536237185, -0.055550357757325114, -0.034402963072744035, -0.02355295203281852, -0.01395989796835577, 0.04791532258364179, 0.005256561332065603, -0.033329144105376754, -0.0469376073953147, -0.03883335799749384, -0.010009045280512503]'
$ws.Range("B7").Value = -0.03838862779801082
$ws.Range("B9").Value = 0.7992911694649822
